# OPPINCOME_holdings.xlsx - refresh the "as of" date in the disclosure note
# and update the model Weight / Percent Change columns with the latest figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet ships protected, so the locked data cells cannot be written
# to until it is unprotected.
$ws.Unprotect()

# Bump the "Model holdings provided as of" date from 2021-03-30 to 2021-03-31.
$ws.Range("A37").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-31 for illustrative purposes only and are subject to change."

# Refresh Weight (D) and Percent Change (E) for every holding row.
$ws.Cells.Item(2, 4).Value = 0.03850304641133154
$ws.Cells.Item(2, 5).Value = 0.007387247278382691

$ws.Cells.Item(3, 4).Value = 0.0216718857810095
$ws.Cells.Item(3, 5).Value = 0.002742946708463867

$ws.Cells.Item(4, 4).Value = 0.02001984216091747
$ws.Cells.Item(4, 5).Value = -0.007646207894192902

$ws.Cells.Item(5, 4).Value = 0.04026679404844045
$ws.Cells.Item(5, 5).Value = -0.004569420035149374

$ws.Cells.Item(6, 4).Value = 0.03744332803936859
$ws.Cells.Item(6, 5).Value = -0.00430191630817367

$ws.Cells.Item(7, 4).Value = 0.02103024092302793
$ws.Cells.Item(7, 5).Value = 0.000776548243059505

$ws.Cells.Item(8, 4).Value = 0.03772487442144042
$ws.Cells.Item(8, 5).Value = 0.008051889957503944

$ws.Cells.Item(9, 4).Value = 0.02130999089546937
$ws.Cells.Item(9, 5).Value = 0.001287474710318204

$ws.Cells.Item(10, 4).Value = 0.0261949186103955
$ws.Cells.Item(10, 5).Value = 0.003013300083125614

$ws.Cells.Item(11, 4).Value = 0.02423792084638112
$ws.Cells.Item(11, 5).Value = -0.002245929253228529

$ws.Cells.Item(12, 4).Value = 0.05854531647430114
$ws.Cells.Item(12, 5).Value = -0.002234359483614723

$ws.Cells.Item(13, 4).Value = 0.02652480474252143
$ws.Cells.Item(13, 5).Value = -0.001847063169560403

$ws.Cells.Item(14, 4).Value = 0.02747026057466055
$ws.Cells.Item(14, 5).Value = 0.01007474813129683

$ws.Cells.Item(15, 4).Value = 0.03572176881024607
$ws.Cells.Item(15, 5).Value = -0.007509813961426803

$ws.Cells.Item(16, 4).Value = 0.01906736400022319
$ws.Cells.Item(16, 5).Value = 0.003323179174743673

$ws.Cells.Item(17, 4).Value = 0.0302341184459968
$ws.Cells.Item(17, 5).Value = -0.01053294922578318

$ws.Cells.Item(18, 4).Value = 0.02402028309782675
$ws.Cells.Item(18, 5).Value = 0.003929727230698177

$ws.Cells.Item(19, 4).Value = 0.1331836325308009
$ws.Cells.Item(19, 5).Value = 0.006711409395973256

$ws.Cells.Item(20, 4).Value = 0.00963572351078322
$ws.Cells.Item(20, 5).Value = 0.001689189189189255

$ws.Cells.Item(21, 4).Value = 0.0159392160305063
$ws.Cells.Item(21, 5).Value = 0.007745822276411074

$ws.Cells.Item(22, 4).Value = 0.01733965342903297
$ws.Cells.Item(22, 5).Value = -0.01526391631531188

$ws.Cells.Item(23, 4).Value = 0.0167683407299139
$ws.Cells.Item(23, 5).Value = -0.02406227883934897

$ws.Cells.Item(24, 4).Value = 0.02161004574039975
$ws.Cells.Item(24, 5).Value = -0.004116117850953116

$ws.Cells.Item(25, 4).Value = 0.01185722898370847
$ws.Cells.Item(25, 5).Value = 0.01753306674869259

$ws.Cells.Item(26, 4).Value = 0.04379973298820494
$ws.Cells.Item(26, 5).Value = -0.003187919463087385

$ws.Cells.Item(27, 4).Value = 0.02552415014173945
$ws.Cells.Item(27, 5).Value = -0.0001962130874131329

$ws.Cells.Item(28, 4).Value = 0.0480729015690005
$ws.Cells.Item(28, 5).Value = 0.003859141341051586

$ws.Cells.Item(29, 4).Value = 0.05717830318223079
$ws.Cells.Item(29, 5).Value = 0.007464079119238631

$ws.Cells.Item(30, 4).Value = 0.01350050830227042
$ws.Cells.Item(30, 5).Value = -0.02476572958500667

$ws.Cells.Item(31, 4).Value = 0.01444612744437595
$ws.Cells.Item(31, 5).Value = 0.002355157795572271

$ws.Cells.Item(32, 4).Value = 0.04432640762984163
$ws.Cells.Item(32, 5).Value = 0.005194805194805197

$ws.Cells.Item(33, 4).Value = 0.01683126950363297
$ws.Cells.Item(33, 5).Value = -0.001584786053882792

# Totals row: Weight stays 100% (1), only Percent Change is refreshed.
$ws.Cells.Item(34, 5).Value = 0.000741046190863015

# Restore protection on the worksheet so it matches the original, shipped
# (protected) state.
$ws.Protect("D382", $true, $true, $true, $false, $false, $false, $false)
